$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 167387
$ws.Range("C4").Value = 158268
$ws.Range("C5").Value = 9120
$ws.Range("C7").Value = 5.45
$ws.Range("C8").Value = 65.37
